# Rebalance NPC Level Stat Table
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Player")

# Row 3 - CLVL1
$ws.Range("C3").Value = 20

# Row 4 - CLVL2
$ws.Range("B4").Value = 130
$ws.Range("C4").Value = 50

# Row 5 - CLVL3
$ws.Range("B5").Value = 160
$ws.Range("C5").Value = 75

# Row 6 - CLVL4
$ws.Range("B6").Value = 190
$ws.Range("C6").Value = 115
$ws.Range("D6").Value = 5

# Row 7 - CLVL5
$ws.Range("B7").Value = 230
$ws.Range("C7").Value = 155
$ws.Range("D7").Value = 5
$ws.Range("G7").Value = 480

# Row 8 - CLVL6
$ws.Range("B8").Value = 330
$ws.Range("C8").Value = 195
$ws.Range("D8").Value = 10
$ws.Range("G8").Value = 480

# Row 9 - CLVL7
$ws.Range("B9").Value = 460
$ws.Range("C9").Value = 245
$ws.Range("D9").Value = 15
$ws.Range("G9").Value = 480

# Row 10 - CLVL8
$ws.Range("B10").Value = 590
$ws.Range("C10").Value = 245
$ws.Range("D10").Value = 15
$ws.Range("G10").Value = 560

# Row 11 - CLVL9
$ws.Range("B11").Value = 720
$ws.Range("D11").Value = 15
$ws.Range("G11").Value = 560

# Row 12 - CLVL10
$ws.Range("B12").Value = 800
$ws.Range("C12").Value = 290
$ws.Range("D12").Value = 20
$ws.Range("G12").Value = 600

# Update the selected cell to H17 on the active sheet view
$ws.Range("H17").Select()
